$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.456.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.689.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "686.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.434"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.312.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.681.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.422.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.650"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.836.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000125"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.64%  "
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.664.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.159"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0898"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "165.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000285"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.22%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("B48").Value = "SuiNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.88%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("E51").Value = "  -2.76%  "
